$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that are present in all three data rows (2, 3, 4) before the edit.
$cols = @("A","B","C","D","E","F","G","H","I","J","K","N","P","Q","R","S","T","U","V","W", `
          "Y","Z","AA","AB","AD","AE","AF","AG","AT","AW","AX","AY")

# Snapshot the current (pre-edit) values of rows 2, 3 and 4 for every shared column,
# plus the value of AC4 (the only one of the three rows that carries a comment),
# before anything gets overwritten.
$row2 = @{}
$row3 = @{}
$row4 = @{}
foreach ($col in $cols) {
    $row2[$col] = $ws.Range("${col}2").Value()
    $row3[$col] = $ws.Range("${col}3").Value()
    $row4[$col] = $ws.Range("${col}4").Value()
}
$ac4 = $ws.Range("AC4").Value()

# Column L only existed (as an empty cell) on row 2. After the rotation below, the
# row-2 record ends up on row 4, so row 4 needs an empty L cell instead. Copying
# the still-untouched K2:L2 pair onto K4:L4 recreates an empty L4 cell (K4 is
# empty on every row anyway, so it is unaffected).
$ws.Range("K2:L2").Copy($ws.Range("K4:L4"))
$ws.Range("L2").Clear()

# The record that used to be on row 2 moves down to row 4, while the records on
# rows 3 and 4 shift up to rows 2 and 3 respectively.
foreach ($col in $cols) {
    $ws.Range("${col}2").Value = $row3[$col]
    $ws.Range("${col}3").Value = $row4[$col]
    $ws.Range("${col}4").Value = $row2[$col]
}

# Column AC only held data (a comment) on row 4; after the shift it belongs on row 3.
$ws.Range("AC3").Value = $ac4
$ws.Range("AC4").Clear()
